$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7 data (date written with a leading apostrophe so it is kept as
# literal text "2024-08-23" instead of being auto-converted to a date
# serial number, then the style is reset back to Normal so no extra
# number-format style gets attached to the cell).
$ws.Range("A7").Value = "'2024-08-23"
$ws.Range("A7").Style = "Normal"

$ws.Range("B7").Value = "II"
$ws.Range("C7").Value = "MOLINO"
$ws.Range("D7").Value = "CABEZAS CONTRERAS KELVIN BRATH"
$ws.Range("E7").Value = "No cumple"

# F7:P7 stay blank (matching the empty cells in the source row), but we
# still want them materialized as real (empty) cells, so briefly flip the
# number format and reset the style back to Normal.
$ws.Range("F7:P7").NumberFormat = "General"
$ws.Range("F7:P7").Style = "Normal"

$ws.Range("Q7").Value = "MAYRA PILAMUNGA"
$ws.Range("R7").Value = "rgrtgr"
